$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.860.92'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.707.54'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.010'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3970'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.47%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4101'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.508'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.007'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.11'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.64%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08959'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.713'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.182'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001355'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.722.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '100.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.41%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07186'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.29%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.407'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.009'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.822.91'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.104'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.344'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.286'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +23.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '166.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '140.10'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.226'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.973'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.44%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09149'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.71%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.088'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.93%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.03079'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +12.84%  '

$ws.Range("E36").Value = '  +3.22%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.976'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.63'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09339'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.49%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.486'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.25%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7842'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.91%  '

$ws.Range("B44").Value = 'NEARProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.645'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.30%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7319'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.17%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.253'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.62%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.005'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.49%  '

$ws.Range("B48").Value = 'Flow'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QQ0NCmjVq+flow-flow'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.356'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.89%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '95.38'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08068'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.10%  '
